$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), copying the formatting used by
# the other header cells (e.g. H1) so the same cell style is reused.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-14
$values = @(
    @(9, 9),
    @(3, 4),
    @(6, 7),
    @(4, 5),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(8, 9),
    @(8, 8),
    @(4, 4),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
